$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Run all the test cases: flip Runmode (column C) from "N" to "Y" for rows 2-17.
$ws.Range("C2:C17").Value = "Y"

# Reflect the new view/selection state captured for this sheet.
$ws.Activate()
[void]$ws.Range("B11").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
